# Update attendee counts / min ticket prices for the "万圣漫控嘉年华10" refresh
# (gh-pages data regenerated at commit 456a3b4).
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 372   # 想去人数: 371 -> 372
$ws1.Range("G2").Value = 70    # 最低票价: 60 -> 70
$ws1.Range("F3").Value = 71    # 想去人数: 70 -> 71
$ws1.Range("F5").Value = 4233  # 想去人数: 4207 -> 4233
$ws1.Range("F7").Value = 457   # 想去人数: 456 -> 457

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 372   # 想去人数: 371 -> 372
$ws4.Range("G2").Value = 70    # 最低票价: 60 -> 70
$ws4.Range("F3").Value = 71    # 想去人数: 70 -> 71
$ws4.Range("F5").Value = 4233  # 想去人数: 4207 -> 4233
$ws4.Range("F9").Value = 457   # 想去人数: 456 -> 457
